# Update "想去人数" (number of people interested) figures in the
# "展览" and "全部类型" sheets to the refreshed values captured at
# commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of old value -> new value for column F, used to locate and update
# the correct rows without depending on row numbers that differ between
# the two sheets.
$updates = @{
    2047  = 2054
    341   = 345
    585   = 590
    95    = 97
    2060  = 2062
    10544 = 10563
    153   = 154
    408   = 410
    7439  = 7459
    1114  = 1115
    709   = 710
    210   = 220
    63    = 64
    3310  = 3315
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.UsedRange
    $lastRow = $usedRange.Rows.Count

    for ($r = 2; $r -le $lastRow; $r++) {
        $cell = $ws.Cells.Item($r, 6)  # column F
        $val = $cell.Value2
        if ($null -ne $val -and $updates.ContainsKey([int]$val)) {
            $cell.Value = $updates[[int]$val]
        }
    }
}
